$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to insert as the new row 2 (B2:F2), shifting existing rows 2-11 down to 3-11
# (the former row 11 data is dropped, as only rows 2-11 exist in the table)
$newRow = @(0.09494813238296555, 0.5332219804412588, 0.4006434568944247, 0.6329640249606803, 0.6398668578903548)

# Capture current values for columns B:F, rows 2 through 10 (these will move to rows 3 through 11)
$colLetters = @("B", "C", "D", "E", "F")
$oldValues = @{}
for ($r = 2; $r -le 10; $r++) {
    $oldValues[$r] = @{}
    foreach ($col in $colLetters) {
        $oldValues[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Shift rows 2-10 down to rows 3-11
for ($r = 10; $r -ge 2; $r--) {
    foreach ($col in $colLetters) {
        $ws.Range("$col$($r+1)").Value = $oldValues[$r][$col]
    }
}

# Write the new row into row 2
for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $ws.Range("$($colLetters[$i])2").Value = $newRow[$i]
}
